# Applies the recomputed Step1/Step2/Step3 values to 34_11R22.xlsx that
# resulted from the "Added Tire Type Filtering for dashboard script and
# cleanup of Tire Type extraction in process_audio_to_csv script" commit.
#
# The external audio-processing pipeline re-derived the per-segment
# Signal_Value_* intensities (Step1_Data). That cascades into the running,
# row-wise cumulative sums (Step2_Sj) and into the derived threshold
# statistics sheets (Step3_DataPts_0.5 / 0.7 / 0.8 / 0.9: first-noticeable-
# increase index/value, threshold-crossing index/value and pulse width).
# This script writes the final values for every cell affected by that
# recomputation.

$wb = $excel.ActiveWorkbook

function Set-Cells($Worksheet, $Values) {
    foreach ($ref in $Values.Keys) {
        $Worksheet.Range($ref).Value = $Values[$ref]
    }
}


# Step1_Data
$ws = $wb.Worksheets.Item("Step1_Data")
$values = @{
    "D2" = 0.03870502709331878
    "E2" = 0.01442294103101579
    "F2" = 0.05481627818721881
    "G2" = 0.127566958481816
    "K2" = 0.01804883216754858
    "M2" = 0.04986564470169137
    "N2" = 0.001962358969917726
    "O2" = 0.1558656093300014
    "R2" = 0.09575161758967751
    "T2" = 0.2568037386478605
    "V2" = 0.08930399832224134
    "X2" = 0.06503752042495724
    "Z2" = 0.002951649283817465
    "AB2" = 0.01722365271253467
    "AC2" = 0.01167417305638264
    "D3" = 0.06030395935270255
    "E3" = 0.01833892476396565
    "F3" = 0.1264042207040559
    "G3" = 0.1009281621099216
    "K3" = 0.006161172916394787
    "M3" = 0.02290097862502213
    "O3" = 0.1387815912424963
    "R3" = 0.03740446883307973
    "T3" = 0.1984627542252333
    "V3" = 0.08984786457455754
    "X3" = 0.0978104781306605
    "Z3" = 0.0204739724653183
    "AB3" = 0.05130134309047839
    "AC3" = 0.02295002717320101
    "AE3" = 0.007930081792912311
    "D4" = 0.05652261696613039
    "E4" = 0.02271715538431906
    "F4" = 0.165753648731275
    "G4" = 0.0956362656625275
    "K4" = 0.02588740307222113
    "M4" = 0.04432237569580003
    "O4" = 0.1749542918933592
    "R4" = 0.03546836936935121
    "T4" = 0.1587505841899322
    "V4" = 0.08371135053624627
    "X4" = 0.0719702556851963
    "Z4" = 0.01057916320763586
    "AB4" = 0.04611142699941614
    "AC4" = 0.007615092606589643
    "E5" = 0.1123365162505709
    "G5" = 0.08023435786795294
    "H5" = 0.0548712682281909
    "I5" = 0.008559740557759649
    "L5" = 0.05051418482231758
    "N5" = 0.04832634546479715
    "O5" = 0.04971525270161967
    "P5" = 0.151924940671481
    "S5" = 0.03894356192885025
    "T5" = 0.0523576996723563
    "U5" = 0.1580301522983452
    "V5" = 0.01205157270082377
    "W5" = 0.04955608083125958
    "X5" = 0.03398299333837992
    "Y5" = 0.04609276335575605
    "Z5" = 0.0001176131962880529
    "AB5" = 0.01340084435790701
    "AC5" = 0.03752828328949263
    "AE5" = 0.001455828465851617
    "E6" = 0.160701448128282
    "F6" = 0.018548893209366
    "G6" = 0.2405999813044859
    "L6" = 0.03463559353624502
    "M6" = 0.03066996784170477
    "O6" = 0.1233111703301069
    "P6" = 0.01090467497612231
    "Q6" = 0.009768924621187827
    "R6" = 0.007150223823547681
    "T6" = 0.1191762215200979
    "U6" = 0.02328755529402197
    "V6" = 0.07460492484545161
    "X6" = 0.07524020089070731
    "Z6" = 0.02241279848880293
    "AB6" = 0.01427928202436894
    "AC6" = 0.02300232741845698
    "AE6" = 0.01170581174704414
}
Set-Cells $ws $values

# Step2_Sj
$ws = $wb.Worksheets.Item("Step2_Sj")
$values = @{
    "D2" = 0.03870502709331878
    "E2" = 0.05312796812433457
    "F2" = 0.1079442463115534
    "G2" = 0.2355112047933694
    "H2" = 0.2355112047933694
    "I2" = 0.2355112047933694
    "J2" = 0.2355112047933694
    "K2" = 0.253560036960918
    "L2" = 0.253560036960918
    "M2" = 0.3034256816626094
    "N2" = 0.3053880406325271
    "O2" = 0.4612536499625285
    "P2" = 0.4612536499625285
    "Q2" = 0.4612536499625285
    "R2" = 0.557005267552206
    "S2" = 0.557005267552206
    "T2" = 0.8138090062000666
    "U2" = 0.8138090062000666
    "V2" = 0.9031130045223079
    "W2" = 0.9031130045223079
    "X2" = 0.9681505249472652
    "Y2" = 0.9681505249472652
    "Z2" = 0.9711021742310826
    "AA2" = 0.9711021742310826
    "AB2" = 0.9883258269436173
    "AC2" = 0.9999999999999999
    "AD2" = 0.9999999999999999
    "AE2" = 0.9999999999999999
    "AF2" = 0.9999999999999999
    "AG2" = 0.9999999999999999
    "AH2" = 0.9999999999999999
    "AI2" = 0.9999999999999999
    "AJ2" = 0.9999999999999999
    "D3" = 0.06030395935270255
    "E3" = 0.0786428841166682
    "F3" = 0.2050471048207241
    "G3" = 0.3059752669306457
    "H3" = 0.3059752669306457
    "I3" = 0.3059752669306457
    "J3" = 0.3059752669306457
    "K3" = 0.3121364398470405
    "L3" = 0.3121364398470405
    "M3" = 0.3350374184720626
    "N3" = 0.3350374184720626
    "O3" = 0.4738190097145589
    "P3" = 0.4738190097145589
    "Q3" = 0.4738190097145589
    "R3" = 0.5112234785476386
    "S3" = 0.5112234785476386
    "T3" = 0.7096862327728719
    "U3" = 0.7096862327728719
    "V3" = 0.7995340973474294
    "W3" = 0.7995340973474294
    "X3" = 0.8973445754780899
    "Y3" = 0.8973445754780899
    "Z3" = 0.9178185479434082
    "AA3" = 0.9178185479434082
    "AB3" = 0.9691198910338866
    "AC3" = 0.9920699182070877
    "AD3" = 0.9920699182070877
    "D4" = 0.05652261696613039
    "E4" = 0.07923977235044945
    "F4" = 0.2449934210817245
    "G4" = 0.340629686744252
    "H4" = 0.340629686744252
    "I4" = 0.340629686744252
    "J4" = 0.340629686744252
    "K4" = 0.3665170898164731
    "L4" = 0.3665170898164731
    "M4" = 0.4108394655122732
    "N4" = 0.4108394655122732
    "O4" = 0.5857937574056324
    "P4" = 0.5857937574056324
    "Q4" = 0.5857937574056324
    "R4" = 0.6212621267749836
    "S4" = 0.6212621267749836
    "T4" = 0.7800127109649158
    "U4" = 0.7800127109649158
    "V4" = 0.863724061501162
    "W4" = 0.863724061501162
    "X4" = 0.9356943171863583
    "Y4" = 0.9356943171863583
    "Z4" = 0.9462734803939942
    "AA4" = 0.9462734803939942
    "AB4" = 0.9923849073934103
    "E5" = 0.1123365162505709
    "F5" = 0.1123365162505709
    "G5" = 0.1925708741185238
    "H5" = 0.2474421423467147
    "I5" = 0.2560018829044744
    "J5" = 0.2560018829044744
    "K5" = 0.2560018829044744
    "L5" = 0.306516067726792
    "M5" = 0.306516067726792
    "N5" = 0.3548424131915891
    "O5" = 0.4045576658932088
    "P5" = 0.5564826065646898
    "Q5" = 0.5564826065646898
    "R5" = 0.5564826065646898
    "S5" = 0.5954261684935401
    "T5" = 0.6477838681658964
    "U5" = 0.8058140204642417
    "V5" = 0.8178655931650655
    "W5" = 0.8674216739963251
    "X5" = 0.901404667334705
    "Y5" = 0.9474974306904611
    "Z5" = 0.9476150438867491
    "AA5" = 0.9476150438867491
    "AB5" = 0.9610158882446561
    "AC5" = 0.9985441715341488
    "AD5" = 0.9985441715341488
    "E6" = 0.160701448128282
    "F6" = 0.179250341337648
    "G6" = 0.419850322642134
    "H6" = 0.419850322642134
    "I6" = 0.419850322642134
    "J6" = 0.419850322642134
    "K6" = 0.419850322642134
    "L6" = 0.454485916178379
    "M6" = 0.4851558840200837
    "N6" = 0.4851558840200837
    "O6" = 0.6084670543501907
    "P6" = 0.619371729326313
    "Q6" = 0.6291406539475008
    "R6" = 0.6362908777710485
    "S6" = 0.6362908777710485
    "T6" = 0.7554670992911464
    "U6" = 0.7787546545851683
    "V6" = 0.8533595794306199
    "W6" = 0.8533595794306199
    "X6" = 0.9285997803213272
    "Y6" = 0.9285997803213272
    "Z6" = 0.9510125788101301
    "AA6" = 0.9510125788101301
    "AB6" = 0.965291860834499
    "AC6" = 0.988294188252956
    "AD6" = 0.988294188252956
}
Set-Cells $ws $values

# Step3_DataPts_0.5
$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$values = @{
    "C2" = 2
    "E2" = 0
    "F2" = 0.557005267552206
    "G2" = 15
    "F3" = 0.5112234785476386
    "C4" = 1
    "E4" = 0
    "F4" = 0.5857937574056324
    "G4" = 13
    "F5" = 0.5564826065646898
    "D6" = 14
    "F6" = 0.6084670543501907
    "G6" = 12
}
Set-Cells $ws $values

# Step3_DataPts_0.7
$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$values = @{
    "C2" = 2
    "E2" = 0
    "F2" = 0.8138090062000666
    "G2" = 17
    "F3" = 0.7096862327728719
    "C4" = 1
    "E4" = 0
    "F4" = 0.7800127109649158
    "G4" = 18
    "F5" = 0.8058140204642417
    "D6" = 19
    "F6" = 0.7554670992911464
    "G6" = 17
}
Set-Cells $ws $values

# Step3_DataPts_0.8
$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$values = @{
    "C2" = 2
    "E2" = 0
    "F2" = 0.8138090062000666
    "G2" = 17
    "D3" = 23
    "F3" = 0.8973445754780899
    "G3" = 22
    "C4" = 1
    "D4" = 21
    "E4" = 0
    "F4" = 0.863724061501162
    "G4" = 20
    "F5" = 0.8058140204642417
    "D6" = 21
    "E6" = 0
    "F6" = 0.8533595794306199
}
Set-Cells $ws $values

# Step3_DataPts_0.9
$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$values = @{
    "C2" = 2
    "E2" = 0
    "F2" = 0.9031130045223079
    "G2" = 19
    "D3" = 25
    "F3" = 0.9178185479434082
    "G3" = 24
    "C4" = 1
    "D4" = 23
    "E4" = 0
    "F4" = 0.9356943171863583
    "G4" = 22
    "D5" = 23
    "F5" = 0.901404667334705
    "G5" = 21
    "D6" = 23
    "E6" = 0
    "F6" = 0.9285997803213272
}
Set-Cells $ws $values
